$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "PO Forecast"

# Header row
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Match the bold / centered / bordered header style used on the other sheets
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4160
$ws.Range("A1:D1").Borders.LineStyle = 1

# Data rows
$ws.Range("A2").Value = 44934.99999999999
$ws.Range("B2").Value = 28
$ws.Range("C2").Value = 9.948845030953201
$ws.Range("D2").Value = 45.82664521648447

$ws.Range("A3").Value = 45004.99999999999
$ws.Range("B3").Value = 269
$ws.Range("C3").Value = 249.3537533395726
$ws.Range("D3").Value = 287.4333241799654

$ws.Range("A4").Value = 45011.99999999999
$ws.Range("B4").Value = 293
$ws.Range("C4").Value = 275.1074702451568
$ws.Range("D4").Value = 311.8442662738559

$ws.Range("A5").Value = 45018.99999999999
$ws.Range("B5").Value = 317
$ws.Range("C5").Value = 299.0577673469685
$ws.Range("D5").Value = 334.2467466495272

$ws.Range("A6").Value = 45025.99999999999
$ws.Range("B6").Value = 341
$ws.Range("C6").Value = 322.7569385994078
$ws.Range("D6").Value = 360.3012810292065

$ws.Range("A7").Value = 45032.99999999999
$ws.Range("B7").Value = 365
$ws.Range("C7").Value = 347.4757180538884
$ws.Range("D7").Value = 385.4175975423831

$ws.Range("A8").Value = 45039.99999999999
$ws.Range("B8").Value = 389
$ws.Range("C8").Value = 370.9866120741506
$ws.Range("D8").Value = 408.101093039058

$ws.Range("A9").Value = 45046.99999999999
$ws.Range("B9").Value = 413
$ws.Range("C9").Value = 394.1145252073844
$ws.Range("D9").Value = 430.8512170204554

$ws.Range("A10").Value = 45053.99999999999
$ws.Range("B10").Value = 437
$ws.Range("C10").Value = 418.8202224387055
$ws.Range("D10").Value = 455.8583164043255

$ws.Range("A11").Value = 45060.99999999999
$ws.Range("B11").Value = 461
$ws.Range("C11").Value = 442.2538526320793
$ws.Range("D11").Value = 478.3852908554921

$ws.Range("A12").Value = 45067.99999999999
$ws.Range("B12").Value = 485
$ws.Range("C12").Value = 466.1295036078751
$ws.Range("D12").Value = 503.4856474333595

# Match the date number format used for the date column on the other sheets
$ws.Range("A2:A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
